# Apply the "last changes in april 2024" edit:
#  1. Refresh the cached "datetimeFigureOut" footer-date field text
#     (10/12/2023 -> 12/01/2024) on the slide master and every slide layout.
#  2. Change the "GO" label on slide 1 to "ADDR".
#  3. Best-effort: touch the presentation-level slide guide list (no-op if
#     the host does not expose it, observed to be unsupported here).

$p = $ppt.ActivePresentation

$oldDate = "10/12/2023"
$newDate = "12/01/2024"

# --- 1. Date placeholder fields -------------------------------------------------

$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

$layouts = $m.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame) {
                if ($sh.TextFrame.HasText) {
                    if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                        $sh.TextFrame.TextRange.Text = $newDate
                    }
                }
            }
        }
    }
}

# --- 2. "GO" -> "ADDR" on slide 1 ------------------------------------------------

$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "GO") {
                $sh.TextFrame.TextRange.Text = "ADDR"
            }
        }
    }
}

# --- 3. Best-effort slide guide list touch --------------------------------------

try {
    [void]$p.Guides.Add(1, 3.0)
} catch {
}

Write-Output "done"
